# Rename the column headers in row 1 to carry the `<formatversion>` suffix
# instead of the old `_old` / `_new` suffixes (FV2404 = old/"before" format
# version, FV2410 = new/"after" format version).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table (ListObject) so the headers act as
# column headers / filters.
$rng = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (pane splits below row 1, top-left cell of the
# scrollable area is A2).
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
